# Add a new worksheet "Sheet2" right after "Sheet1" and populate it,
# making it the active sheet (mirrors the target workbook/sheet diff).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = 123
$ws2.Range("A2").Value = "test"

$ws2.Range("A3").Select() | Out-Null
